$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 166, shifting existing rows 166-232 down to 167-233
$ws.Rows.Item(166).Insert()

# Populate the new row 166 with data
$ws.Cells.Item(166, 1).Value = 10
$ws.Cells.Item(166, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(166, 3).Value = "La Araucanía"
$ws.Cells.Item(166, 4).Value = 44704
$ws.Cells.Item(166, 4).NumberFormat = $ws.Cells.Item(167, 4).NumberFormat
$ws.Cells.Item(166, 5).Value = 9
$ws.Cells.Item(166, 6).Value = "Fruta"
$ws.Cells.Item(166, 7).Value = 100103
$ws.Cells.Item(166, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(166, 9).Value = 100103002
$ws.Cells.Item(166, 10).Value = "Ciruela"
$ws.Cells.Item(166, 11).Value = "Pink Delight"
$ws.Cells.Item(166, 12).Value = "Primera"
$ws.Cells.Item(166, 13).Value = 5
$ws.Cells.Item(166, 14).Value = 260000
$ws.Cells.Item(166, 15).Value = 260000
$ws.Cells.Item(166, 16).Value = 260000
$ws.Cells.Item(166, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(166, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(166, 19).Value = 578
$ws.Cells.Item(166, 20).Value = 450
